$d = $word.ActiveDocument

# 1. The npm-install-dependencies line had its text split across two runs
#    ("...express-validator" + " "). Re-issuing a Find/Replace over the
#    whole phrase (now including the trailing space) collapses it back
#    into a single run while keeping the original run formatting.
$d.Content.Find.Execute(
    "npm i express express-handlebars express-session mysql express-mysql-session morgan bcryptjs passport passport-local timeago.js connect-flash express-validator ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "npm i express express-handlebars express-session mysql express-mysql-session morgan bcryptjs passport passport-local timeago.js connect-flash express-validator ",
    2) | Out-Null

# 2. The "nodemon --save-dev" explanation sentence was split across four
#    runs (" npm i " + "--save-dev " + "nodemon" + ", el -d es..."). A
#    Find/Replace across the full sentence merges them into one run.
$d.Content.Find.Execute(
    " npm i --save-dev nodemon, el -d es para indicar que es una dependencia de desarrollo",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " npm i --save-dev nodemon, el -d es para indicar que es una dependencia de desarrollo",
    2) | Out-Null

# 3. The "se ejecuta el script npm run dev" line was split across two
#    runs ("-se ejecuta el script " + "npm run dev"); merge them too.
$d.Content.Find.Execute(
    "-se ejecuta el script npm run dev",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "-se ejecuta el script npm run dev",
    2) | Out-Null

# 4. Append two new notes right after the paragraph that explains how
#    index.js is created in the routes folder, reusing that paragraph's
#    formatting (Lucida Console, sz 18).
$anchorText = "-se crea index.js en la carpeta routes para configurar las rutas de navegacion"
$anchor = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd("`r", "`n")
    if ($t -eq $anchorText) {
        $anchor = $p
        break
    }
}

if ($anchor -ne $null) {
    $anchor.Range.InsertParagraphAfter()

    $newP1 = $anchor.Next()
    $newP1.Range.Text = "-para un textarea no existe la propiedad value, si se quiere asignar un valor a un textarea, tiene que ser entre los caracteres de mayor y menor de la etiqueta."

    $newP1.Range.InsertParagraphAfter()
    $newP2 = $newP1.Next()
    $newP2.Range.Text = "-nodemon reinicia servidor cuando se modifica y graba un javascript, cuando se modifica html, no lo reinicia automaticamente."
}
